# Fruta / hortaliza, semanal
# Weekly update: insert a new price record as row 26 (Berenjena,
# Terminal Hortofrutícola Agro Chillán), pushing the existing rows
# 26-30 down to 27-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 26, shifting rows 26:30
# down to 27:31.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with this week's record.
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44694
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112001
$ws.Range("G26").Value = "Berenjena"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 10000
$ws.Range("N26").Value = "$/caja 60 unidades"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 167
$ws.Range("Q26").Value = 60
$ws.Range("R26").Value = "Hortaliza"
